# Insert two new price records for "Feria Lagunitas de Puerto Montt - Lechuga"
# as new rows 434 and 435, pushing the existing rows (old 434..459) down to
# become rows 436..461 (the last existing row is effectively duplicated at
# the very end of the block as row 461, matching the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 434 downwards by inserting two new blank rows above the
# existing row 434, copying that row's formatting (this preserves the
# YYYY-MM-DD HH:MM:SS date style on column D for the new rows too).
$ws.Rows("434:435").Insert()

# --- New row 434 ---
$ws.Cells.Item(434, 1).Value = 4
$ws.Cells.Item(434, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(434, 3).Value = "Los Lagos"
$ws.Cells.Item(434, 4).Value = [DateTime]"2022-02-18"
$ws.Cells.Item(434, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(434, 5).Value = 10
$ws.Cells.Item(434, 6).Value = 100112033
$ws.Cells.Item(434, 7).Value = "Lechuga"
$ws.Cells.Item(434, 8).Value = "Escarola"
$ws.Cells.Item(434, 9).Value = "Primera"
$ws.Cells.Item(434, 10).Value = 250
$ws.Cells.Item(434, 11).Value = 12000
$ws.Cells.Item(434, 12).Value = 12000
$ws.Cells.Item(434, 13).Value = 12000
$ws.Cells.Item(434, 14).Value = "$/caja 15 unidades"
$ws.Cells.Item(434, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(434, 16).Value = 800
$ws.Cells.Item(434, 17).Value = 15
$ws.Cells.Item(434, 18).Value = "Hortaliza"

# --- New row 435 ---
$ws.Cells.Item(435, 1).Value = 4
$ws.Cells.Item(435, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(435, 3).Value = "Los Lagos"
$ws.Cells.Item(435, 4).Value = [DateTime]"2022-02-18"
$ws.Cells.Item(435, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(435, 5).Value = 10
$ws.Cells.Item(435, 6).Value = 100112033
$ws.Cells.Item(435, 7).Value = "Lechuga"
$ws.Cells.Item(435, 8).Value = "Escarola"
$ws.Cells.Item(435, 9).Value = "Segunda"
$ws.Cells.Item(435, 10).Value = 250
$ws.Cells.Item(435, 11).Value = 10000
$ws.Cells.Item(435, 12).Value = 10000
$ws.Cells.Item(435, 13).Value = 10000
$ws.Cells.Item(435, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(435, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(435, 16).Value = 556
$ws.Cells.Item(435, 17).Value = 18
$ws.Cells.Item(435, 18).Value = "Hortaliza"
